# Update ticket-interest / price data for the "展览" and "全部类型" sheets.
# Both sheets contain the same 16 events (rows 2-17) that need refreshed
# "想去人数" (F) figures, and row 5's "最低票价" (G) flips from a numeric
# price to the text "已售罄" (sold out) because F5 (interest count) also grew.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 127
    $ws.Range("F3").Value = 406
    $ws.Range("F4").Value = 11914

    $ws.Range("F5").Value = 1257
    $ws.Range("G5").Value = "已售罄"

    $ws.Range("F7").Value = 23
    $ws.Range("F10").Value = 179
    $ws.Range("F11").Value = 180
    $ws.Range("F13").Value = 57
    $ws.Range("F15").Value = 36
    $ws.Range("F17").Value = 1532
}
